# Team Roles.xlsx update
# - Marks the "Tai" analysis rows for McKinley(2), Taft(1), Harding(1),
#   Hoover(1), F.Roosevelt(3) and Truman(2) as "In Progress" (was "-"),
#   matching the blue "In Progress" styling already used elsewhere in
#   the Status column.
# - Scrolls the sheet view back to the top and updates the active
#   selection to G8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Status cells to flip from "-" to "In Progress"
$rows = @(2, 5, 8, 11, 14, 17)

foreach ($row in $rows) {
    $cell = $ws.Range("F" + $row)
    $cell.Value = "In Progress"
    $cell.Font.Color = 12611584   # RGB(0, 112, 192) == the existing "In Progress" blue
}

# Reset the view: scroll back to the top-left (A1) and select G8
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G8").Select()
